$wb = $excel.ActiveWorkbook

# --- Sheet "Q.No.5": add the tokenizer instructions + link in C2:C3 ---
$ws5 = $wb.Worksheets.Item("Q.No.5")
$ws5.Range("C2").Value = "Use the tokenizer to calculate the tokens for different prompts (with space ,with special characters without space etc)"
$ws5.Range("C3").Value = "https://platform.openai.com/tokenizer"
[void]$ws5.Range("C2:C3").Select()

# --- Sheet "Q.No.6": insert 3 new rows above the existing table and add the new instruction lines ---
$ws6 = $wb.Worksheets.Item("Q.No.6")
$ws6.Rows("3:5").Insert()
$ws6.Range("A2").Value = "From Openrouter ai --Extract the values of "
$ws6.Range("A3").Value = "1)Context window"
$ws6.Range("A4").Value = "2)Input tokens(with price)"
$ws6.Range("A5").Value = "3)Output token (with price)"
[void]$ws6.Range("F4").Select()
